# The deck currently applies the "Integral" (Red Violet) design theme
# (ppt/theme/theme2.xml, the theme wired to the slide master / presentation
# design) while ppt/theme/theme1.xml (wired only to the notes master)
# carries the stock "Office Theme" palette.
#
# The target edit swaps the two palettes: the design actually used across
# the slides switches from the Red-Violet "Integral" colors to the
# standard "Office Theme" colors (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# PowerPoint exposes the active design's 12 theme colors through
# ThemeColorScheme (Slide/Master/Layout all resolve to the same design),
# indexed 1-12 in the fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hlink, 12 folHlink
# RGB values are written using the standard VBA RGB() packing:
#   RGB(r,g,b) = r + g*256 + b*65536

$p = $ppt.ActivePresentation

function LongRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the stock "Office Theme" color scheme.
$officeTheme = @{
    1  = (LongRGB 0x00 0x00 0x00)   # dk1      000000
    2  = (LongRGB 0xFF 0xFF 0xFF)   # lt1      FFFFFF
    3  = (LongRGB 0x44 0x54 0x6A)   # dk2      44546A
    4  = (LongRGB 0xE7 0xE6 0xE6)   # lt2      E7E6E6
    5  = (LongRGB 0x5B 0x9B 0xD5)   # accent1  5B9BD5
    6  = (LongRGB 0xED 0x7D 0x31)   # accent2  ED7D31
    7  = (LongRGB 0xA5 0xA5 0xA5)   # accent3  A5A5A5
    8  = (LongRGB 0xFF 0xC0 0x00)   # accent4  FFC000
    9  = (LongRGB 0x44 0x72 0xC4)   # accent5  4472C4
    10 = (LongRGB 0x70 0xAD 0x47)   # accent6  70AD47
    11 = (LongRGB 0x05 0x63 0xC1)   # hlink    0563C1
    12 = (LongRGB 0x95 0x4F 0x72)   # folHlink 954F72
}

# Apply through the first slide's ThemeColorScheme - it addresses the
# presentation's single active design theme (shared by the slide master,
# every layout/slide and the notes master alike).
$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeTheme[$i]
}
